# Actualización automática de tasas-transfi.xlsx
# Updates the daily conversion-rate figures on "Hoja1" (summary text in A1)
# and the underlying rate cells on "tasas" (N10/O10/N12/O12).

$wb = $excel.ActiveWorkbook

# --- Hoja1: A1 summary text -------------------------------------------------
$ws1 = $wb.Worksheets.Item("Hoja1")

$oldText = $ws1.Range("A1").Value2
$newText = $oldText.Replace(
    "1000 Bs = 9.92 = 42252.96 pesos",
    "1000 Bs = 9.76 = 41512.08 pesos"
).Replace(
    "42252.96 pesos = 9.91 = 972.1 Bs",
    "41512.08 pesos = 9.71 = 951.62 Bs"
)
$ws1.Range("A1").Value = $newText

# --- tasas: rate cells -------------------------------------------------------
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 102.455
$ws2.Range("O10").Value = 4253.12
$ws2.Range("N12").Value = 4274.99
$ws2.Range("O12").Value = 98
